# UPDATE[calc]: Range operator now working properly (#260)
#
# The CHOOSE function's range operator support is now working, so the two
# example cells that previously held these formulas as literal (quoted)
# text can now contain real, calculating formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B10 used to be the literal text "=SUM(A9:CHOOSE(2,A10,A11,A12))" and C10
# held a "not supported yet" note next to it. Turn B10 into a real formula
# and clear the now-unnecessary note in C10.
$ws.Range("B10").Formula = "=SUM(A9:CHOOSE(2,A10,A11,A12))"
$ws.Range("C10").ClearContents()

# B11 used to be the literal text "=SUM(CHOOSE(2,A9:A10,A9:A11,A9:A12))".
# Turn it into a real, calculating formula as well.
$ws.Range("B11").Formula = "=SUM(CHOOSE(2,A9:A10,A9:A11,A9:A12))"

# Reflect the author's selection at save time.
$null = $ws.Range("B11").Select()
